# Apply the table-style change made to the "SOURCES OF FINANCE" table on
# slide 6 (the table was re-styled from the deck's custom "Table_0" style
# to the built-in PowerPoint table style {FFAA6DBE-1B31-4A5A-810B-D5C5AD9CDF1F},
# e.g. via the Table Design ribbon's style gallery).

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(6)

# Locate the shape that actually hosts the table (it's the second shape on
# this slide - the first is the slide title textbox).
for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $sh = $s.Shapes.Item($i)
    if ($sh.HasTable) {
        $sh.Table.ApplyStyle("{FFAA6DBE-1B31-4A5A-810B-D5C5AD9CDF1F}")
    }
}
